# Auto-generated Excel COM-interop edit script
# Applies the Bahamut_Profits.xlsx numeric updates (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row 86
$ws_ALC.Range("H86").Value = 66670280
$ws_ALC.Range("I86").Value = 5951.5
$ws_ALC.Range("J86").Value = 76926330
$ws_ALC.Range("K86").Value = 5951.5
$ws_ALC.Range("L86").Value = 76926330
$ws_ALC.Range("M86").Value = -4828.5
$ws_ALC.Range("N86").Value = -76928576

# ALC!row 89
$ws_ALC.Range("H89").Value = 66670280
$ws_ALC.Range("I89").Value = 5951.5
$ws_ALC.Range("J89").Value = 76926330
$ws_ALC.Range("K89").Value = 29757.5
$ws_ALC.Range("L89").Value = 384631650
$ws_ALC.Range("M89").Value = -24141.5
$ws_ALC.Range("N89").Value = -384642882

# ALC!row 92
$ws_ALC.Range("H92").Value = 2736.2856
$ws_ALC.Range("I92").Value = 2944.8235
$ws_ALC.Range("J92").Value = 1850
$ws_ALC.Range("K92").Value = 2944.8235
$ws_ALC.Range("L92").Value = 1850
$ws_ALC.Range("M92").Value = -1696.8235
$ws_ALC.Range("N92").Value = -4346

# ALC!row 94
$ws_ALC.Range("H94").Value = 665
$ws_ALC.Range("I94").Value = 665
$ws_ALC.Range("K94").Value = 665
$ws_ALC.Range("M94").Value = -214

# ALC!row 103
$ws_ALC.Range("H103").Value = 1118.3158
$ws_ALC.Range("I103").Value = 1330.7273
$ws_ALC.Range("J103").Value = 826.25
$ws_ALC.Range("K103").Value = 3992.1819
$ws_ALC.Range("L103").Value = 2478.75
$ws_ALC.Range("M103").Value = -3406.1819
$ws_ALC.Range("N103").Value = -3650.75

# ALC!row 106
$ws_ALC.Range("H106").Value = 3612.1428
$ws_ALC.Range("I106").Value = 3580.8333
$ws_ALC.Range("K106").Value = 3580.8333
$ws_ALC.Range("M106").Value = -2949.8333

# ALC!row 125
$ws_ALC.Range("H125").Value = 29558
$ws_ALC.Range("I125").Value = 24258
$ws_ALC.Range("J125").Value = 33798
$ws_ALC.Range("K125").Value = 218322
$ws_ALC.Range("L125").Value = 304182
$ws_ALC.Range("M125").Value = -215862
$ws_ALC.Range("N125").Value = -309102

# ALC!row 132
$ws_ALC.Range("H132").Value = 1540.1702
$ws_ALC.Range("I132").Value = 1541.381
$ws_ALC.Range("K132").Value = 4624.143
$ws_ALC.Range("M132").Value = -2094.143

# ALC!row 137
$ws_ALC.Range("H137").Value = 1100.4857
$ws_ALC.Range("I137").Value = 909.4286
$ws_ALC.Range("J137").Value = 1864.7142
$ws_ALC.Range("K137").Value = 2728.2858
$ws_ALC.Range("L137").Value = 5594.142599999999
$ws_ALC.Range("M137").Value = -178.2857999999997
$ws_ALC.Range("N137").Value = -10694.1426

# ALC!row 138
$ws_ALC.Range("H138").Value = 4201.727
$ws_ALC.Range("I138").Value = 885.0526
$ws_ALC.Range("J138").Value = 4989.4375
$ws_ALC.Range("K138").Value = 2655.1578
$ws_ALC.Range("L138").Value = 14968.3125
$ws_ALC.Range("M138").Value = 2484.8422
$ws_ALC.Range("N138").Value = -25248.3125

# ARM!row 74
$ws_ARM.Range("H74").Value = 814.05
$ws_ARM.Range("I74").Value = 814.7895
$ws_ARM.Range("K74").Value = 814.7895
$ws_ARM.Range("M74").Value = 59.21050000000002

# ARM!row 77
$ws_ARM.Range("H77").Value = 814.05
$ws_ARM.Range("I77").Value = 814.7895
$ws_ARM.Range("K77").Value = 4073.9475
$ws_ARM.Range("M77").Value = 294.0525000000002

# ARM!row 110
$ws_ARM.Range("H110").Value = 547.9545000000001
$ws_ARM.Range("I110").Value = 547.9545000000001
$ws_ARM.Range("J110").Value = 0
$ws_ARM.Range("K110").Value = 547.9545000000001
$ws_ARM.Range("L110").Value = 0
$ws_ARM.Range("M110").Value = 1497.0455
$ws_ARM.Range("N110").ClearContents()

# ARM!row 112
$ws_ARM.Range("H112").Value = 2517500
$ws_ARM.Range("J112").Value = 2517500
$ws_ARM.Range("L112").Value = 2517500
$ws_ARM.Range("N112").Value = -2520454

# ARM!row 132
$ws_ARM.Range("H132").Value = 1337.625
$ws_ARM.Range("I132").Value = 577.13336
$ws_ARM.Range("J132").Value = 2605.111
$ws_ARM.Range("K132").Value = 1731.40008
$ws_ARM.Range("L132").Value = 7815.333
$ws_ARM.Range("M132").Value = 798.5999199999999
$ws_ARM.Range("N132").Value = -12875.333

# BSM!row 134
$ws_BSM.Range("H134").Value = 44782.5
$ws_BSM.Range("I134").Value = 3734.9333
$ws_BSM.Range("K134").Value = 11204.7999
$ws_BSM.Range("M134").Value = -8669.7999

# CRP!row 2
$ws_CRP.Range("H2").Value = 500
$ws_CRP.Range("J2").Value = 500
$ws_CRP.Range("L2").Value = 500
$ws_CRP.Range("N2").Value = -726

# CRP!row 31
$ws_CRP.Range("H31").Value = 2436.8635
$ws_CRP.Range("I31").Value = 2489.1667
$ws_CRP.Range("J31").Value = 2201.5
$ws_CRP.Range("K31").Value = 2489.1667
$ws_CRP.Range("L31").Value = 2201.5
$ws_CRP.Range("M31").Value = -2194.1667
$ws_CRP.Range("N31").Value = -2791.5

# CRP!row 34
$ws_CRP.Range("H34").Value = 2436.8635
$ws_CRP.Range("I34").Value = 2489.1667
$ws_CRP.Range("J34").Value = 2201.5
$ws_CRP.Range("K34").Value = 2489.1667
$ws_CRP.Range("L34").Value = 2201.5
$ws_CRP.Range("M34").Value = -2287.1667
$ws_CRP.Range("N34").Value = -2605.5

# CRP!row 58
$ws_CRP.Range("H58").Value = 4352.393
$ws_CRP.Range("I58").Value = 616.8261
$ws_CRP.Range("J58").Value = 21536
$ws_CRP.Range("K58").Value = 616.8261
$ws_CRP.Range("L58").Value = 21536
$ws_CRP.Range("M58").Value = -413.8261
$ws_CRP.Range("N58").Value = -21942

# CRP!row 104
$ws_CRP.Range("H104").Value = 24856.666
$ws_CRP.Range("J104").Value = 24856.666
$ws_CRP.Range("L104").Value = 24856.666
$ws_CRP.Range("N104").Value = -30098.666

# CRP!row 136
$ws_CRP.Range("H136").Value = 4352.393
$ws_CRP.Range("I136").Value = 616.8261
$ws_CRP.Range("J136").Value = 21536
$ws_CRP.Range("K136").Value = 1850.4783
$ws_CRP.Range("L136").Value = 64608
$ws_CRP.Range("M136").Value = 699.5217
$ws_CRP.Range("N136").Value = -69708

# CUL!row 48
$ws_CUL.Range("H48").Value = 0
$ws_CUL.Range("I48").Value = 0
$ws_CUL.Range("K48").Value = 0
$ws_CUL.Range("M48").ClearContents()

# CUL!row 97
$ws_CUL.Range("H97").Value = 11904878
$ws_CUL.Range("I97").Value = 23809690
$ws_CUL.Range("J97").Value = 67.333336
$ws_CUL.Range("K97").Value = 71429070
$ws_CUL.Range("L97").Value = 202.000008
$ws_CUL.Range("M97").Value = -71428574
$ws_CUL.Range("N97").Value = -1194.000008

# CUL!row 98
$ws_CUL.Range("H98").Value = 2711.111
$ws_CUL.Range("I98").Value = 490
$ws_CUL.Range("J98").Value = 3821.6667
$ws_CUL.Range("K98").Value = 1470
$ws_CUL.Range("L98").Value = 11465.0001
$ws_CUL.Range("M98").Value = 28
$ws_CUL.Range("N98").Value = -14461.0001

# CUL!row 101
$ws_CUL.Range("H101").Value = 4875
$ws_CUL.Range("J101").Value = 5428.5713
$ws_CUL.Range("L101").Value = 16285.7139
$ws_CUL.Range("N101").Value = -21153.7139

# CUL!row 129
$ws_CUL.Range("H129").Value = 92814.09
$ws_CUL.Range("J129").Value = 169010.83
$ws_CUL.Range("L129").Value = 507032.49
$ws_CUL.Range("N129").Value = -517032.49

# CUL!row 131
$ws_CUL.Range("H131").Value = 12694.262
$ws_CUL.Range("J131").Value = 1471.3768
$ws_CUL.Range("L131").Value = 4414.1304
$ws_CUL.Range("N131").Value = -14494.1304

# CUL!row 132
$ws_CUL.Range("H132").Value = 1083.1666
$ws_CUL.Range("J132").Value = 1083.1666
$ws_CUL.Range("L132").Value = 9748.499400000001
$ws_CUL.Range("N132").Value = -14808.4994

# CUL!row 137
$ws_CUL.Range("H137").Value = 4537.8
$ws_CUL.Range("I137").Value = 1498.5714
$ws_CUL.Range("J137").Value = 6174.3076
$ws_CUL.Range("K137").Value = 4495.7142
$ws_CUL.Range("L137").Value = 18522.9228
$ws_CUL.Range("M137").Value = 604.2857999999997
$ws_CUL.Range("N137").Value = -28722.9228

# GSM!row 2
$ws_GSM.Range("H2").Value = 116.333336
$ws_GSM.Range("I2").Value = 90.833336
$ws_GSM.Range("J2").Value = 167.33333
$ws_GSM.Range("K2").Value = 90.833336
$ws_GSM.Range("L2").Value = 167.33333
$ws_GSM.Range("M2").Value = 22.166664
$ws_GSM.Range("N2").Value = -393.33333

# GSM!row 48
$ws_GSM.Range("H48").Value = 5000
$ws_GSM.Range("I48").Value = 5000
$ws_GSM.Range("J48").Value = 0
$ws_GSM.Range("K48").Value = 5000
$ws_GSM.Range("L48").Value = 0
$ws_GSM.Range("M48").Value = -4515
$ws_GSM.Range("N48").ClearContents()

# GSM!row 132
$ws_GSM.Range("H132").Value = 3833.375
$ws_GSM.Range("I132").Value = 3536.3333
$ws_GSM.Range("J132").Value = 4724.5
$ws_GSM.Range("K132").Value = 10608.9999
$ws_GSM.Range("L132").Value = 14173.5
$ws_GSM.Range("M132").Value = -8078.999899999999
$ws_GSM.Range("N132").Value = -19233.5

# LTW!row 16
$ws_LTW.Range("H16").Value = 1000
$ws_LTW.Range("I16").Value = 1000
$ws_LTW.Range("J16").Value = 1000
$ws_LTW.Range("K16").Value = 1000
$ws_LTW.Range("L16").Value = 1000
$ws_LTW.Range("M16").Value = -830
$ws_LTW.Range("N16").Value = -1340

# WVR!row 130
$ws_WVR.Range("H130").Value = 29980
$ws_WVR.Range("J130").Value = 29980
$ws_WVR.Range("L130").Value = 29980
$ws_WVR.Range("N130").Value = -40020

# WVR!row 132
$ws_WVR.Range("H132").Value = 2693.125
$ws_WVR.Range("I132").Value = 2148.7
$ws_WVR.Range("K132").Value = 6446.099999999999
$ws_WVR.Range("M132").Value = -3916.099999999999

# WVR!row 136
$ws_WVR.Range("H136").Value = 1834.1333
$ws_WVR.Range("I136").Value = 1950.0312
$ws_WVR.Range("K136").Value = 5850.0936
$ws_WVR.Range("M136").Value = -3300.0936
